# Update the API table:
# - Row 7 ("process" request) gains K7="final.jpg" and the L7 description is
#   extended to also cover sending the final.jpg file.
# - The old row 8 ("get_pic" for "final") is removed entirely; it is merged
#   into row 7's description.
# - The old row 9 ("end") shifts up to become row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing K value for the "process" row and extend its description.
$ws.Range("K7").Value = "final.jpg"
$ws.Range("L7").Value = "initiate processing, save as final.jpg to output directory and send file final.jpg from output directory"

# Delete the old row 8 (get_pic / final row) - this shifts row 9 ("end") up to row 8.
$ws.Rows(8).Delete()

# The "index" column holds literal sequential numbers, not formulas, so
# deleting a row does not renumber it automatically - fix the "end" row's
# index (was 6, should now be 5 since a row was removed).
$ws.Range("A8").Value = 5

# Cosmetic: reflect that the user ended up with the (now last) row selected
# and zoomed in, matching the saved view state after this edit.
$ws.Application.ActiveWindow.Zoom = 157
[void]$ws.Rows(8).Select()
